# Update Metrics sheet source values (B2:B13). Downstream sheets (e.g. "today")
# pull these via formulas (Metrics!B2 ... Metrics!B13) and will recalculate
# automatically.
$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 271358.64000000007
$metrics.Range("B3").Value = 239126.15000000002
$metrics.Range("B4").Value = 83809.110000000015
$metrics.Range("B5").Value = 11051
$metrics.Range("B6").Value = 5067604.3900000015
$metrics.Range("B7").Value = 4281202.83
$metrics.Range("B8").Value = 1490768.9400000002
$metrics.Range("B9").Value = 197258
$metrics.Range("B10").Value = 33532985.38000001
$metrics.Range("B11").Value = 31556477.990000002
$metrics.Range("B12").Value = 11772490.98
$metrics.Range("B13").Value = 1294888

# Move the saved selection on the Metrics sheet from D22 to E22.
$metrics.Range("E22").Select()

# Move the saved selection on the "today" sheet (the active/visible sheet)
# from E8 to F8.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F8").Select()
